$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Clear old data in the table region that will be rewritten
$ws.Range("A2:E9").ClearContents()

# Row 2 (summary row) - B2 and D2 text stay same, just shared string index shifts naturally
$ws.Range("B2").Value = "duracion total"
$ws.Range("C2").Value = 73229.78606
$ws.Range("D2").Value = "No representativa porque fue en paralelo!"

# Row 5 - header row
$ws.Range("A5").Value = "dataset"
$ws.Range("B5").Value = "tiempo promedio por configuración"
$ws.Range("C5").Value = "error mínimo obtenido"
$ws.Range("D5").Value = "error promedio"
$ws.Range("E5").Value = "H"
$ws.Range("F5").Value = "alpha"
$ws.Range("G5").Value = "'nu"

# Row 6 - complex
$ws.Range("A6").Value = "complex"
$ws.Range("B6").Value = 17.1
$ws.Range("C6").Value = 0.089109
$ws.Range("D6").Value = 0.3644
$ws.Range("E6").Value = 20
$ws.Range("F6").Value = 0.0000005
$ws.Range("G6").Value = 0.01

# Row 7 - linear
$ws.Range("A7").Value = "linear"
$ws.Range("B7").Value = 12.9
$ws.Range("C7").Value = 0
$ws.Range("D7").Value = 0.1967
$ws.Range("E7").Value = 3
$ws.Range("F7").Value = 0.0000001
$ws.Range("G7").Value = 0.005

# Row 8 - ring
$ws.Range("A8").Value = "ring"
$ws.Range("B8").Value = 12.9
$ws.Range("C8").Value = 0.008
$ws.Range("D8").Value = 0.1261
$ws.Range("E8").Value = 3
$ws.Range("F8").Value = 0.0000001
$ws.Range("G8").Value = 0.01

# Row 9 - xor
$ws.Range("A9").Value = "xor"
$ws.Range("B9").Value = 6.1
$ws.Range("C9").Value = 0
$ws.Range("D9").Value = 0.0459
$ws.Range("E9").Value = 5
$ws.Range("F9").Value = 0.0000001
$ws.Range("G9").Value = 0.005

# Number format for alpha column (scientific notation)
$ws.Range("F6:F9").NumberFormat = "0.00E+00"

# Column C was manually resized (no longer auto bestFit)
$ws.Range("C1").ColumnWidth = 17.8

$ws.Range("A5:G9").Select()
